$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    21  = "Name and contact information if an agent is being used."
    25  = "Name and contact information if an agent is being used."
    33  = "Telephone number and email address of the applicant."
    37  = "Name and contact information for the parties making the application."
    43  = "How any natural habitats on the development site will be improved by the proposed works."
    58  = "Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation."
    59  = "Details of any conflict of interest that may exist between the applicant and planning authority."
    62  = "Signed and dated verification of the application's accuracy."
    65  = "How the proposed development will impact employment, including existing and proposed employee numbers"
    72  = "How the site is currently being used."
    82  = "Results of any flood risk assessments made for the development site"
    88  = "Proposed operating hours if the proposed development is intended for non-residential use."
    95  = "Details of changes to non-residential floorspace in the proposed development."
    112 = "Who will be affected by the proposal and whether they have been notified, such as agricultural tenants"
    127 = "Details of pre-application advice received from the planning authority"
    132 = "How waste will be managed on the site"
    142 = "What development, works or change of use is proposed"
    148 = "Details of the residential units that make up both the current and proposed development."
    164 = "How big the site is including relevant measurements"
    166 = "Where the proposed development will be built."
    175 = "Information to help the planning authority arrange a site visit"
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
